# Wallet workbook update: add support for american stocks too.
# - Adds BBAS3.SA in A5 (new Brazilian ticker)
# - Adds AAPL in A6 (new American ticker)
# - Moves the active selection to A7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tickers added to the wallet list
$ws.Range("A5").Value2 = "BBAS3.SA"
$ws.Range("A6").Value2 = "AAPL"

# Leave the selection where the user would continue typing the next ticker
$ws.Range("A7").Select()
